$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry data to be rotated among rows 14, 15, 16.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Capture original values (rows 14, 15, 16) before overwriting anything.
$row14 = @{}
$row15 = @{}
$row16 = @{}
foreach ($c in $cols) {
    $row14[$c] = $ws.Range($c + "14").Value2
    $row15[$c] = $ws.Range($c + "15").Value2
    $row16[$c] = $ws.Range($c + "16").Value2
}

# Apply rotation:
#   new row14 = old row15
#   new row15 = old row16
#   new row16 = old row14
foreach ($c in $cols) {
    $ws.Range($c + "14").Value = $row15[$c]
    $ws.Range($c + "15").Value = $row16[$c]
    $ws.Range($c + "16").Value = $row14[$c]
}
